$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to literal text (prevents Excel auto-converting numeric-looking
# strings such as "1.011" or "0.000008811" into actual numbers), matching the source
# data which stores every Price/Volume cell as inline text.
function Set-TextValue($cell, [string]$text) {
    if ($text -match "^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$") {
        $cell.Value = "'" + $text
    } else {
        $cell.Value = $text
    }
}

# Row 2
Set-TextValue $ws.Range("D2") "27.428.19"
Set-TextValue $ws.Range("E2") "  +1.59%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.860.08"
Set-TextValue $ws.Range("E3") "  +0.62%  "

# Row 4
Set-TextValue $ws.Range("E4") "  -0.05%  "

# Row 5
Set-TextValue $ws.Range("D5") "311.48"
Set-TextValue $ws.Range("E5") "  +0.76%  "

# Row 6
Set-TextValue $ws.Range("E6") "  +0.03%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.4768"
Set-TextValue $ws.Range("E7") "  -0.16%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.3810"
Set-TextValue $ws.Range("E8") "  +3.62%  "

# Row 9
Set-TextValue $ws.Range("E9") "  +1.17%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.9302"
Set-TextValue $ws.Range("E10") "  -0.07%  "

# Row 11
Set-TextValue $ws.Range("D11") "20.80"
Set-TextValue $ws.Range("E11") "  +5.08%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.07786"
Set-TextValue $ws.Range("E12") "  +0.67%  "

# Row 13
Set-TextValue $ws.Range("D13") "1.861.37"
Set-TextValue $ws.Range("E13") "  -0.09%  "

# Row 14
Set-TextValue $ws.Range("D14") "5.438"
Set-TextValue $ws.Range("E14") "  +1.72%  "

# Row 15
Set-TextValue $ws.Range("E15") "  +1.55%  "

# Row 16
Set-TextValue $ws.Range("D16") "90.09"
Set-TextValue $ws.Range("E16") "  +1.13%  "

# Row 17
Set-TextValue $ws.Range("D17") "1.011"
Set-TextValue $ws.Range("E17") "  -0.14%  "

# Row 18
Set-TextValue $ws.Range("D18") "0.000008811"
Set-TextValue $ws.Range("E18") "  +2.05%  "

# Row 19
Set-TextValue $ws.Range("E19") "  -0.12%  "

# Row 20
Set-TextValue $ws.Range("D20") "27.391.33"
Set-TextValue $ws.Range("E20") "  +1.33%  "

# Row 21
Set-TextValue $ws.Range("D21") "14.63"
Set-TextValue $ws.Range("E21") "  +0.77%  "

# Row 22
Set-TextValue $ws.Range("D22") "5.093"
Set-TextValue $ws.Range("E22") "  +0.45%  "

# Row 23
Set-TextValue $ws.Range("E23") "  +0.45%  "

# Row 24
Set-TextValue $ws.Range("D24") "1.941"
Set-TextValue $ws.Range("E24") "  +0.57%  "

# Row 25
Set-TextValue $ws.Range("D25") "155.57"
Set-TextValue $ws.Range("E25") "  +1.85%  "

# Row 26
Set-TextValue $ws.Range("E26") "  +1.32%  "

# Row 27
Set-TextValue $ws.Range("E27") "  -0.29%  "

# Row 28
Set-TextValue $ws.Range("D28") "115.23"
Set-TextValue $ws.Range("E28") "  +0.90%  "

# Row 29
Set-TextValue $ws.Range("D29") "4.951"
Set-TextValue $ws.Range("E29") "  -0.22%  "

# Row 30
Set-TextValue $ws.Range("D30") "0.08897"
Set-TextValue $ws.Range("E30") "  +0.51%  "

# Row 31
Set-TextValue $ws.Range("E31") "  +0.45%  "

# Row 32
Set-TextValue $ws.Range("D32") "1.205"
Set-TextValue $ws.Range("E32") "  +2.12%  "

# Row 33
Set-TextValue $ws.Range("D33") "4.592"
Set-TextValue $ws.Range("E33") "  +1.98%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.7504"
Set-TextValue $ws.Range("E34") "  +1.42%  "

# Row 35
Set-TextValue $ws.Range("D35") "2.721"
Set-TextValue $ws.Range("E35") "  -0.27%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.02044"
Set-TextValue $ws.Range("E36") "  +4.23%  "

# Row 37
Set-TextValue $ws.Range("D37") "1.124"
Set-TextValue $ws.Range("E37") "  +1.34%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.5538"
Set-TextValue $ws.Range("E38") "  +5.41%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.05263"
Set-TextValue $ws.Range("E39") "  -0.02%  "

# Row 40
Set-TextValue $ws.Range("D40") "2.987"
Set-TextValue $ws.Range("E40") "  +0.45%  "

# Row 41
Set-TextValue $ws.Range("D41") "7.028"
Set-TextValue $ws.Range("E41") "  +0.16%  "

# Row 42
Set-TextValue $ws.Range("D42") "8.603"
Set-TextValue $ws.Range("E42") "  +4.32%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.1523"
Set-TextValue $ws.Range("E43") "  +0.45%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D44") "10.69"
Set-TextValue $ws.Range("E44") "  +0.54%  "

# Row 45
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D45") "0.4870"
Set-TextValue $ws.Range("E45") "  +2.54%  "

# Row 46
Set-TextValue $ws.Range("E46") "  -0.06%  "

# Row 47
Set-TextValue $ws.Range("D47") "1.665"
Set-TextValue $ws.Range("E47") "  +3.77%  "

# Row 48
Set-TextValue $ws.Range("D48") "102.75"
Set-TextValue $ws.Range("E48") "  +0.96%  "

# Row 49
Set-TextValue $ws.Range("D49") "67.28"
Set-TextValue $ws.Range("E49") "  +2.34%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.06092"
Set-TextValue $ws.Range("E50") "  +0.41%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.9120"
Set-TextValue $ws.Range("E51") "  +2.58%  "
